$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted before the existing row 687,
# shifting all subsequent rows (687-701) down by one (to 688-702).
$ws.Rows.Item(687).Insert()

$ws.Cells.Item(687, 1).Value = 6
$ws.Cells.Item(687, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(687, 3).Value = "Metropolitana"
$ws.Cells.Item(687, 4).Value = 45239
$ws.Cells.Item(687, 5).Value = 13
$ws.Cells.Item(687, 6).Value = 100112043
$ws.Cells.Item(687, 7).Value = "Pepino ensalada"
$ws.Cells.Item(687, 8).Value = "Sin especificar"
$ws.Cells.Item(687, 9).Value = "Primera"
$ws.Cells.Item(687, 10).Value = 590
$ws.Cells.Item(687, 11).Value = 15000
$ws.Cells.Item(687, 12).Value = 16000
$ws.Cells.Item(687, 13).Value = 15458
$ws.Cells.Item(687, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(687, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(687, 16).Value = 258
$ws.Cells.Item(687, 17).Value = 60
$ws.Cells.Item(687, 18).Value = "Hortaliza"
